$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 is new; copy the date style (s="2") used by column A in existing rows
# so the new date cell A53 is formatted consistently with A2:A52.
$ws.Range("A52").Copy($ws.Range("A53"))

$ws.Range("A2").Value2 = 39400
$ws.Range("B2").Value2 = 2007
$ws.Range("C2").Value2 = 4.930115226412335
$ws.Range("D2").Value2 = 2008
$ws.Range("E2").Value2 = 1.402909115503936

$ws.Range("A3").Value2 = 39583
$ws.Range("B3").Value2 = 2008
$ws.Range("C3").Value2 = 2.319409367208825
$ws.Range("D3").Value2 = 2009
$ws.Range("E3").Value2 = 5.080273296954374

$ws.Range("A4").Value2 = 39765
$ws.Range("B4").Value2 = 2008
$ws.Range("C4").Value2 = 1.457587285166628
$ws.Range("D4").Value2 = 2009
$ws.Range("E4").Value2 = 2.807231216534278

$ws.Range("A5").Value2 = 39948
$ws.Range("B5").Value2 = 2009
$ws.Range("C5").Value2 = -3.942037578692481
$ws.Range("D5").Value2 = 2010
$ws.Range("E5").Value2 = -1.648748515828502

$ws.Range("A6").Value2 = 40130
$ws.Range("B6").Value2 = 2009
$ws.Range("C6").Value2 = -0.9140166223623569
$ws.Range("D6").Value2 = 2010
$ws.Range("E6").Value2 = 1.821983295885121

$ws.Range("A7").Value2 = 40310
$ws.Range("B7").Value2 = 2010
$ws.Range("C7").Value2 = -2.839753013810475
$ws.Range("D7").Value2 = 2011
$ws.Range("E7").Value2 = -1.632723506456935

$ws.Range("A8").Value2 = 40494
$ws.Range("B8").Value2 = 2010
$ws.Range("C8").Value2 = 2.585942866987878
$ws.Range("D8").Value2 = 2011
$ws.Range("E8").Value2 = 2.722861752007866

$ws.Range("A9").Value2 = 40676
$ws.Range("B9").Value2 = 2011
$ws.Range("C9").Value2 = 4.960109259035406
$ws.Range("D9").Value2 = 2012
$ws.Range("E9").Value2 = 4.506881698240095

$ws.Range("A10").Value2 = 40862
$ws.Range("B10").Value2 = 2011
$ws.Range("C10").Value2 = 4.253963781362402
$ws.Range("D10").Value2 = 2012
$ws.Range("E10").Value2 = 2.878414118480799

$ws.Range("A11").Value2 = 41044
$ws.Range("B11").Value2 = 2012
$ws.Range("C11").Value2 = 2.134646894829806
$ws.Range("D11").Value2 = 2013
$ws.Range("E11").Value2 = 2.531943146540772

$ws.Range("A12").Value2 = 41228
$ws.Range("B12").Value2 = 2012
$ws.Range("C12").Value2 = 1.752870900283909
$ws.Range("D12").Value2 = 2013
$ws.Range("E12").Value2 = 3.144721336271927

$ws.Range("A13").Value2 = 41409
$ws.Range("B13").Value2 = 2013
$ws.Range("C13").Value2 = -2.700325749999488
$ws.Range("D13").Value2 = 2014
$ws.Range("E13").Value2 = -0.3858735870725938

$ws.Range("A14").Value2 = 41592
$ws.Range("B14").Value2 = 2013
$ws.Range("C14").Value2 = -1.479696720105139
$ws.Range("D14").Value2 = 2014
$ws.Range("E14").Value2 = 2.238623952069552

$ws.Range("A15").Value2 = 41774
$ws.Range("B15").Value2 = 2014
$ws.Range("C15").Value2 = 5.469647210234996
$ws.Range("D15").Value2 = 2015
$ws.Range("E15").Value2 = 3.061326532789543

$ws.Range("A16").Value2 = 41957
$ws.Range("B16").Value2 = 2014
$ws.Range("C16").Value2 = 3.900127535411246
$ws.Range("D16").Value2 = 2015
$ws.Range("E16").Value2 = -0.6155071485167585

$ws.Range("A17").Value2 = 42137
$ws.Range("B17").Value2 = 2015
$ws.Range("C17").Value2 = 0.950153436409007
$ws.Range("D17").Value2 = 2016
$ws.Range("E17").Value2 = 2.074800935750787

$ws.Range("A18").Value2 = 42321
$ws.Range("B18").Value2 = 2015
$ws.Range("C18").Value2 = 0.03947433952959933
$ws.Range("D18").Value2 = 2016
$ws.Range("E18").Value2 = 1.182212550358064

$ws.Range("A19").Value2 = 42503
$ws.Range("B19").Value2 = 2016
$ws.Range("C19").Value2 = 3.458696398997052
$ws.Range("D19").Value2 = 2017
$ws.Range("E19").Value2 = 2.610227683091337

$ws.Range("A20").Value2 = 42689
$ws.Range("B20").Value2 = 2016
$ws.Range("C20").Value2 = 2.192778679161944
$ws.Range("D20").Value2 = 2017
$ws.Range("E20").Value2 = -0.5835597102573198

$ws.Range("A21").Value2 = 42867
$ws.Range("B21").Value2 = 2017
$ws.Range("C21").Value2 = 2.772413308959698
$ws.Range("D21").Value2 = 2018
$ws.Range("E21").Value2 = 2.755099409670958

$ws.Range("A22").Value2 = 43053
$ws.Range("B22").Value2 = 2017
$ws.Range("C22").Value2 = 3.40836448860673
$ws.Range("D22").Value2 = 2018
$ws.Range("E22").Value2 = 2.34299484087257

$ws.Range("A23").Value2 = 43145
$ws.Range("B23").Value2 = 2018
$ws.Range("C23").Value2 = 1.376958470962375
$ws.Range("D23").Value2 = 2019
$ws.Range("E23").Value2 = 2.554746761760573

$ws.Range("A24").Value2 = 43235
$ws.Range("B24").Value2 = 2018
$ws.Range("C24").Value2 = 3.145819842658448
$ws.Range("D24").Value2 = 2019
$ws.Range("E24").Value2 = 3.875106770584158

$ws.Range("A25").Value2 = 43326
$ws.Range("B25").Value2 = 2018
$ws.Range("C25").Value2 = 2.594480907596441
$ws.Range("D25").Value2 = 2019
$ws.Range("E25").Value2 = 2.978856953395592

$ws.Range("A26").Value2 = 43418
$ws.Range("B26").Value2 = 2018
$ws.Range("C26").Value2 = 2.799070570134488
$ws.Range("D26").Value2 = 2019
$ws.Range("E26").Value2 = 3.056075254339996

$ws.Range("A27").Value2 = 43510
$ws.Range("B27").Value2 = 2019
$ws.Range("C27").Value2 = 3.665688413913726
$ws.Range("D27").Value2 = 2020
$ws.Range("E27").Value2 = 3.126710782028086

$ws.Range("A28").Value2 = 43600
$ws.Range("B28").Value2 = 2019
$ws.Range("C28").Value2 = 4.520465362328063
$ws.Range("D28").Value2 = 2020
$ws.Range("E28").Value2 = 4.124307769579505

$ws.Range("A29").Value2 = 43691
$ws.Range("B29").Value2 = 2019
$ws.Range("C29").Value2 = 4.143226503463837
$ws.Range("D29").Value2 = 2020
$ws.Range("E29").Value2 = 2.872814227731846

$ws.Range("A30").Value2 = 43783
$ws.Range("B30").Value2 = 2019
$ws.Range("C30").Value2 = 4.195393191694419
$ws.Range("D30").Value2 = 2020
$ws.Range("E30").Value2 = 2.652948310315506

$ws.Range("A31").Value2 = 43875
$ws.Range("B31").Value2 = 2020
$ws.Range("C31").Value2 = 2.017172445889925
$ws.Range("D31").Value2 = 2021
$ws.Range("E31").Value2 = 2.549404806734712

$ws.Range("A32").Value2 = 43966
$ws.Range("B32").Value2 = 2020
$ws.Range("C32").Value2 = 4.479055418855871
$ws.Range("D32").Value2 = 2021
$ws.Range("E32").Value2 = 4.888255652935936

$ws.Range("A33").Value2 = 44068
$ws.Range("B33").Value2 = 2020
$ws.Range("C33").Value2 = 2.133862376612439
$ws.Range("D33").Value2 = 2021
$ws.Range("E33").Value2 = 0.273062749117492

$ws.Range("A34").Value2 = 44159
$ws.Range("B34").Value2 = 2020
$ws.Range("C34").Value2 = 1.666553973046048
$ws.Range("D34").Value2 = 2021
$ws.Range("E34").Value2 = -0.9999522486825452

$ws.Range("A35").Value2 = 44251
$ws.Range("B35").Value2 = 2021
$ws.Range("C35").Value2 = -2.991914460776179
$ws.Range("D35").Value2 = 2022
$ws.Range("E35").Value2 = -0.7897024150736676

$ws.Range("A36").Value2 = 44341
$ws.Range("B36").Value2 = 2021
$ws.Range("C36").Value2 = 2.922623512367206
$ws.Range("D36").Value2 = 2022
$ws.Range("E36").Value2 = 2.761298099516418

$ws.Range("A37").Value2 = 44432
$ws.Range("B37").Value2 = 2021
$ws.Range("C37").Value2 = 1.773820722495745
$ws.Range("D37").Value2 = 2022
$ws.Range("E37").Value2 = 1.707328096006822

$ws.Range("A38").Value2 = 44525
$ws.Range("B38").Value2 = 2021
$ws.Range("C38").Value2 = 1.879266440112803
$ws.Range("D38").Value2 = 2022
$ws.Range("E38").Value2 = -0.1343977949472275

$ws.Range("A39").Value2 = 44617
$ws.Range("B39").Value2 = 2022
$ws.Range("C39").Value2 = -0.05874878927536464
$ws.Range("D39").Value2 = 2023
$ws.Range("E39").Value2 = 1.647170385765651

$ws.Range("A40").Value2 = 44706
$ws.Range("B40").Value2 = 2022
$ws.Range("C40").Value2 = -0.2939924376199055
$ws.Range("D40").Value2 = 2023
$ws.Range("E40").Value2 = 1.109922826784815

$ws.Range("A41").Value2 = 44798
$ws.Range("B41").Value2 = 2022
$ws.Range("C41").Value2 = -2.404913754290983
$ws.Range("D41").Value2 = 2023
$ws.Range("E41").Value2 = -2.467583266971451

$ws.Range("A42").Value2 = 44890
$ws.Range("B42").Value2 = 2022
$ws.Range("C42").Value2 = -2.620683231370946
$ws.Range("D42").Value2 = 2023
$ws.Range("E42").Value2 = -2.83913279674276

$ws.Range("A43").Value2 = 44981
$ws.Range("B43").Value2 = 2023
$ws.Range("C43").Value2 = -3.938880009788048
$ws.Range("D43").Value2 = 2024
$ws.Range("E43").Value2 = -0.6223811084860253

$ws.Range("A44").Value2 = 45071
$ws.Range("B44").Value2 = 2023
$ws.Range("C44").Value2 = -2.429883624035745
$ws.Range("D44").Value2 = 2024
$ws.Range("E44").Value2 = -0.8739913853863412

$ws.Range("A45").Value2 = 45163
$ws.Range("B45").Value2 = 2023
$ws.Range("C45").Value2 = -2.901570548279864
$ws.Range("D45").Value2 = 2024
$ws.Range("E45").Value2 = -1.444783760697688

$ws.Range("A46").Value2 = 45254
$ws.Range("B46").Value2 = 2023
$ws.Range("C46").Value2 = -3.036556262700274
$ws.Range("D46").Value2 = 2024
$ws.Range("E46").Value2 = -1.803491225663911

$ws.Range("A47").Value2 = 45345
$ws.Range("B47").Value2 = 2024
$ws.Range("C47").Value2 = -2.786207191231715
$ws.Range("D47").Value2 = 2025
$ws.Range("E47").Value2 = -2.423357831382655

$ws.Range("A48").Value2 = 45436
$ws.Range("B48").Value2 = 2024
$ws.Range("C48").Value2 = -1.060428249734879
$ws.Range("D48").Value2 = 2025
$ws.Range("E48").Value2 = -1.285332554730334

$ws.Range("A49").Value2 = 45534
$ws.Range("B49").Value2 = 2024
$ws.Range("C49").Value2 = -2.859191689251428
$ws.Range("D49").Value2 = 2025
$ws.Range("E49").Value2 = -3.927752965551978

$ws.Range("A50").Value2 = 45618
$ws.Range("B50").Value2 = 2024
$ws.Range("C50").Value2 = -2.953443685011514
$ws.Range("D50").Value2 = 2025
$ws.Range("E50").Value2 = -2.75492543068685

$ws.Range("A51").Value2 = 45713
$ws.Range("B51").Value2 = 2025
$ws.Range("C51").Value2 = -0.8799002368073738
$ws.Range("D51").Value2 = 2026
$ws.Range("E51").Value2 = -2.26429946512674

$ws.Range("A52").Value2 = 45800
$ws.Range("B52").Value2 = 2025
$ws.Range("C52").Value2 = 0.1432745330888219
$ws.Range("D52").Value2 = 2026
$ws.Range("E52").Value2 = -0.484131235569496

$ws.Range("A53").Value2 = 45891
$ws.Range("B53").Value2 = 2025
$ws.Range("C53").Value2 = -1.131442475565558
$ws.Range("D53").Value2 = 2026
$ws.Range("E53").Value2 = -2.567308877837815

Write-Host "Done updating rows 2-53"